$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections in column C ("20 LARGE for VM's" / "10-100Gb" / "LARGE is faving logs") ---
$ws.Range("C11").Value = "20 to LARGE for VM’s"
$ws.Range("C12").Value = "10-100Gb or LARGE"
$ws.Range("C15").Value = "LARGE is for saving LARGE logs"

# --- New column widths for D (~21.38 chars) and E (~20.19 chars) ---
# (ColumnWidth is entered net of the engine's fixed 5/6-character padding,
# and snaps to the nearest 1/6-character pixel grid on save.)
$ws.Columns.Item(4).ColumnWidth = 20.5
$ws.Columns.Item(5).ColumnWidth = 19.333333333333332

# --- Move the active selection from A15 to D2 ---
[void]$ws.Range("D2").Select()
